# Updated cryptos list on Sat Jul 13 23:52:43 UTC 2024 with GitHub Actions
# Refresh prices / 1h volume % for each coin row, and fix the ordering of
# the WrappedEther / Polkadot rows (17 <-> 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.202.21'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '3.175.45'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '532.16'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +11.18%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.27'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('E10').Value = '  +6.15%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.112'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.49%  '
$ws.Range('D12').Value = '3.729.21'
$ws.Range('E12').Value = '  +1.71%  '
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.89'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000171'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.95%  '
$ws.Range('D16').Value = '59.244.86'
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.24'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.123.77'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.03'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '8.17'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '376.14'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('E23').Value = '  +5.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '69.74'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.45'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +15.48%  '
$ws.Range('D28').Value = '0.0₃0872'
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '22.45'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.61%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.04'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.21'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('E34').Value = '  +4.17%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '157.18'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.34%  '
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0715'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +6.30%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '25.40'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.49%  '
$ws.Range('D39').Value = '2.705.16'
$ws.Range('E39').Value = '  +7.52%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.69'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('E41').Value = '  +3.79%  '
$ws.Range('E42').Value = '  +8.55%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.725'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.56%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '39.15'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.50%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Value = '3.221.23'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('E47').Value = '  +12.36%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.986'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.21'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '20.27'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.759'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.77%  '
